$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "8_281115_1914_7_xgboost_with_3in1_preprocess_valid1_valid2_"
$ws.Range("B9").Value = 0.622
$ws.Range("C9").Value = "ensembled 7 xgboost, in 3in1 data set with features preprocessed, with 2 valid sets"

$ws.Range("C10").Select()
